$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 0.27

$ws.Range("B3").Value = 2.78
$ws.Range("C3").Value = 0.26
$ws.Range("D3").Value = 1.32

$ws.Range("B4").Value = 2.13
$ws.Range("C4").Value = 3
